# Applies the "Updated cryptos list" price/volume refresh (and the
# NEARProtocol/ApeXProtocol row-46/47 swap) described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '45.894.31'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.446.32'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.63%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.73%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.89'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.518'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.74%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.531'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.75'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0803'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.35%  '
$ws.Range("E12").Value = '  -1.72%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.14'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.04'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.24%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.831.19'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.72%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.443.65'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.837'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '45.749.81'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.48'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.88%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.39'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0931'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.24'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '245.98'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.34'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.77%  '
$ws.Range("E25").Value = '  +0.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '25.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.71%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.19'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.56%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.68'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.64'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '49.25'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.53%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.129'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.95%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.85'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.90%  '
$ws.Range("E34").Value = '  +2.57%  '
$ws.Range("E35").Value = '  +0.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0760'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.52'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.35%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.89'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.92'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '125.87'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.26'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.36%  '
$ws.Range("E42").Value = '  +1.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.98'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0292'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.955.21'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.39%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.12'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.90%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.95'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.14%  '
$ws.Range("E48").Value = '  +9.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.09'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.99%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '77.43'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.83%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.91'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.72%  '
